$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2252.8157
$ws.Range("J17").Value = 2252.8157
$ws.Range("L17").Value = 6758.4471
$ws.Range("N17").Value = -7094.4471
$ws.Range("H98").Value = 1150.6666
$ws.Range("I98").Value = 990.6667
$ws.Range("J98").Value = 1470.6666
$ws.Range("K98").Value = 990.6667
$ws.Range("L98").Value = 1470.6666
$ws.Range("M98").Value = 507.3333
$ws.Range("N98").Value = -4466.6666
$ws.Range("H122").Value = 1150.6666
$ws.Range("I122").Value = 990.6667
$ws.Range("J122").Value = 1470.6666
$ws.Range("K122").Value = 2972.0001
$ws.Range("L122").Value = 4411.9998
$ws.Range("M122").Value = -522.0001000000002
$ws.Range("N122").Value = -9311.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1097.8235
$ws.Range("I97").Value = 859.4545000000001
$ws.Range("K97").Value = 859.4545000000001
$ws.Range("M97").Value = -363.4545000000001
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3019.2
$ws.Range("I86").Value = 1769
$ws.Range("J86").Value = 5675.875
$ws.Range("K86").Value = 1769
$ws.Range("L86").Value = 5675.875
$ws.Range("M86").Value = -646
$ws.Range("N86").Value = -7921.875
$ws.Range("H89").Value = 3019.2
$ws.Range("I89").Value = 1769
$ws.Range("J89").Value = 5675.875
$ws.Range("K89").Value = 8845
$ws.Range("L89").Value = 28379.375
$ws.Range("M89").Value = -3229
$ws.Range("N89").Value = -39611.375
$ws.Range("I94").Value = 517.25
$ws.Range("J94").Value = 624.25
$ws.Range("K94").Value = 517.25
$ws.Range("L94").Value = 624.25
$ws.Range("M94").Value = -66.25
$ws.Range("N94").Value = -1526.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 72.31579000000001
$ws.Range("I7").Value = 50
$ws.Range("K7").Value = 50
$ws.Range("M7").Value = 63
$ws.Range("H22").Value = 2259.75
$ws.Range("J22").Value = 4750
$ws.Range("L22").Value = 4750
$ws.Range("N22").Value = -5450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 125052.875
$ws.Range("I4").Value = 250016.5
$ws.Range("J4").Value = 89.25
$ws.Range("K4").Value = 750049.5
$ws.Range("L4").Value = 267.75
$ws.Range("M4").Value = -749937.5
$ws.Range("N4").Value = -491.75
$ws.Range("H97").Value = 569.1111
$ws.Range("J97").Value = 397.5
$ws.Range("L97").Value = 1192.5
$ws.Range("N97").Value = -2184.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("K5").Value = 500
$ws.Range("M5").Value = -388
$ws.Range("H9").Value = 167.1
$ws.Range("I9").Value = 212.66667
$ws.Range("J9").Value = 98.75
$ws.Range("K9").Value = 212.66667
$ws.Range("L9").Value = 98.75
$ws.Range("M9").Value = -42.66667000000001
$ws.Range("N9").Value = -438.75
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("N28").Value = ""
$ws.Range("H41").Value = 2359.6
$ws.Range("I41").Value = 2413.6667
$ws.Range("J41").Value = 2278.5
$ws.Range("K41").Value = 2413.6667
$ws.Range("L41").Value = 2278.5
$ws.Range("N41").Value = -2988.5
$ws.Range("H44").Value = 4000
$ws.Range("J44").Value = 4000
$ws.Range("N44").Value = -5192
$ws.Range("H47").Value = 14013.5
$ws.Range("I47").Value = 12028
$ws.Range("J47").Value = 15999
$ws.Range("K47").Value = 12028
$ws.Range("M47").Value = -11460
$ws.Range("N47").Value = -17135
$ws.Range("H49").Value = 10000
$ws.Range("J49").Value = 10000
$ws.Range("N49").Value = -10368
$ws.Range("H97").Value = 1224.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1224.5
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = -2216.5
$ws.Range("H134").Value = 95636.625
$ws.Range("J134").Value = 95636.625
$ws.Range("L134").Value = 286909.875
$ws.Range("N134").Value = -291979.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 909.375
$ws.Range("I22").Value = 912.5
$ws.Range("K22").Value = 912.5
$ws.Range("M22").Value = -617.5
$ws.Range("H27").Value = 909.375
$ws.Range("I27").Value = 912.5
$ws.Range("K27").Value = 912.5
$ws.Range("M27").Value = -805.5
$ws.Range("H68").Value = 6642.6924
$ws.Range("I68").Value = 3971
$ws.Range("J68").Value = 8312.5
$ws.Range("K68").Value = 3971
$ws.Range("L68").Value = 8312.5
$ws.Range("M68").Value = -3222
$ws.Range("N68").Value = -9810.5
$ws.Range("H71").Value = 6642.6924
$ws.Range("I71").Value = 3971
$ws.Range("J71").Value = 8312.5
$ws.Range("K71").Value = 19855
$ws.Range("L71").Value = 41562.5
$ws.Range("M71").Value = -16111
$ws.Range("N71").Value = -49050.5
$ws.Range("H93").Value = 2921.5
$ws.Range("I93").Value = 3843
$ws.Range("K93").Value = 3843
$ws.Range("M93").Value = -2595
$ws.Range("H122").Value = 3079.4
$ws.Range("I122").Value = 2999
$ws.Range("J122").Value = 3200
$ws.Range("K122").Value = 8997
$ws.Range("L122").Value = 9600
$ws.Range("M122").Value = -6547
$ws.Range("N122").Value = -14500

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 16180.667
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 5000
$ws.Range("M69").Value = -4251
$ws.Range("H72").Value = 16180.667
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 15000
$ws.Range("M72").Value = -11256
$ws.Range("H107").Value = 1199.4286
$ws.Range("I107").Value = 932.3333
$ws.Range("K107").Value = 2796.9999
$ws.Range("M107").Value = -876.9998999999998
$ws.Range("H122").Value = 913.84
$ws.Range("I122").Value = 888.4545000000001
$ws.Range("K122").Value = 2665.3635
$ws.Range("M122").Value = -215.3635000000004
